$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '61.929.91'
$c.ClearFormats()
$ws.Range("E2").Value = '  +0.36%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.917.68'
$c.ClearFormats()
$ws.Range("E3").Value = '  -0.42%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range("E4").Value = '  +0.02%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '588.06'
$c.ClearFormats()
$ws.Range("E5").Value = '  -1.37%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '146.30'
$c.ClearFormats()
$ws.Range("E6").Value = '  +3.00%  '
$ws.Range("E7").Value = '  +0.20%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.506'
$c.ClearFormats()
$ws.Range("E8").Value = '  +1.36%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '2.917.10'
$c.ClearFormats()
$ws.Range("E9").Value = '  -0.44%  '
$ws.Range("E10").Value = '  -2.00%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.151'
$c.ClearFormats()
$ws.Range("E11").Value = '  +6.95%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.436'
$c.ClearFormats()
$ws.Range("E12").Value = '  -1.60%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.0000237'
$c.ClearFormats()
$ws.Range("E13").Value = '  +6.68%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '32.46'
$c.ClearFormats()
$ws.Range("E14").Value = '  -1.99%  '
$ws.Range("E15").Value = '  -1.27%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '3.402.21'
$c.ClearFormats()
$ws.Range("E16").Value = '  -0.18%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '61.940.82'
$c.ClearFormats()
$ws.Range("E17").Value = '  +0.69%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '6.61'
$c.ClearFormats()
$ws.Range("E18").Value = '  -0.65%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '2.915.78'
$c.ClearFormats()
$ws.Range("E19").Value = '  -0.63%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '434.60'
$c.ClearFormats()
$ws.Range("E20").Value = '  +0.20%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '13.41'
$c.ClearFormats()
$ws.Range("E21").Value = '  -0.62%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.661'
$c.ClearFormats()
$ws.Range("E22").Value = '  -1.42%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '6.94'
$c.ClearFormats()
$ws.Range("E23").Value = '  -1.50%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '80.97'
$c.ClearFormats()
$ws.Range("E24").Value = '  -0.12%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '10.99'
$c.ClearFormats()
$ws.Range("E25").Value = '  +2.88%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '11.89'
$c.ClearFormats()
$ws.Range("E26").Value = '  +1.87%  '
$ws.Range("E27").Value = '  -1.03%  '
$ws.Range("E28").Value = '  -0.12%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '7.34'
$c.ClearFormats()
$ws.Range("E29").Value = '  +7.27%  '
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.0000104'
$c.ClearFormats()
$ws.Range("E30").Value = '  +20.94%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '2.58'
$c.ClearFormats()
$ws.Range("E31").Value = '  -0.45%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '2.12'
$c.ClearFormats()
$ws.Range("E32").Value = '  +1.50%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.109'
$c.ClearFormats()
$ws.Range("E33").Value = '  +3.45%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '26.09'
$c.ClearFormats()
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("E35").Value = '  +0.13%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.978'
$c.ClearFormats()
$ws.Range("E36").Value = '  -0.58%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '3.08'
$c.ClearFormats()
$ws.Range("E37").Value = '  +8.52%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '5.54'
$c.ClearFormats()
$ws.Range("E38").Value = '  -0.49%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '49.23'
$c.ClearFormats()
$ws.Range("E39").Value = '  +0.02%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '2.00'
$c.ClearFormats()
$ws.Range("E40").Value = '  +2.87%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '8.38'
$c.ClearFormats()
$ws.Range("E41").Value = '  -1.29%  '
$ws.Range("E42").Value = '  -1.70%  '
$ws.Range("E43").Value = '  +0.29%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '39.06'
$c.ClearFormats()
$ws.Range("E44").Value = '  +1.24%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '2.698.14'
$c.ClearFormats()
$ws.Range("E45").Value = '  +0.43%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '135.05'
$c.ClearFormats()
$ws.Range("E46").Value = '  +1.24%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.0338'
$c.ClearFormats()
$ws.Range("E47").Value = '  +0.44%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '348.04'
$c.ClearFormats()
$ws.Range("E48").Value = '  -2.44%  '
$ws.Range("E50").Value = '  +0.52%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '22.57'
$c.ClearFormats()
$ws.Range("E51").Value = '  -0.58%  '
